# "finished running all experiments"
# Row 1538's ID was a zero-padded text placeholder ("00002598"); now that
# the run for that particle finished, it is recorded as the real numeric
# ID (2598) and 38 more freshly-finished experiment rows (IDs 2599-2636)
# are appended below it. The new last row (1577, ID 2637) is still the
# in-flight placeholder, so it keeps the zero-padded text form, exactly
# like row 1538 used to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the now-finished row 1538: was text "00002598", becomes numeric 2598
$ws.Cells.Item(1538, 1).Value = 2598

# H-column (num_particles) results for the newly finished experiments,
# rows 1539 (ID 2599) through 1577 (ID 2637), in order.
$hValues = @(4118, 4024, 4043, 4027, 4103, 4037, 3947, 3946, 3891, 4030, `
             4038, 4101, 4086, 3986, 3936, 3971, 4007, 4080, 3918, 4012, `
             4015, 3971, 3990, 3951, 4020, 4053, 3940, 3984, 3961, 3996, `
             3909, 4038, 4055, 3981, 4066, 3960, 3924, 4055, 3947)

$startRow = 1539
$startId = 2599
$lastRow = 1577

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $startRow + $i
    $id = $startId + $i

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = 0.4
    $ws.Cells.Item($row, 3).Value = 0.0001
    $ws.Cells.Item($row, 4).Value = 0.5
    $ws.Cells.Item($row, 5).Value = 0.8
    $ws.Cells.Item($row, 6).Value = 100
    $ws.Cells.Item($row, 7).Value = 100
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
    $ws.Cells.Item($row, 9).Value = 5000000
}

# The still-running experiment's row (the new last row) keeps its ID as
# zero-padded text, matching how row 1538 looked before it finished.
$ws.Cells.Item($lastRow, 1).Value = "'00002637"
$ws.Cells.Item($lastRow, 1).Style = "Normal"
